$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This applies a roster reshuffle within the Philadelphia 76ers roster table.
# Three players (P.J. Tucker, Shake Milton, Tobias Harris) rotate through rows 4-6,
# and three players (Montrezl Harrell, James Harden, Tyrese Maxey) rotate through rows 9-11.
# All other rows / cells are left untouched.

# Row 4: Shake Milton
$ws.Range("B4").Value = 18
$ws.Range("C4").Value = "Shake Milton"
$ws.Range("D4").Value = "SG"
$ws.Range("E4").Value = "6-5"
$ws.Range("F4").Value = 205
$ws.Range("G4").Value = "September 26, 1996"
$ws.Range("H4").Value = "us"
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = "SMU"
$ws.Range("K4").Value = "https://www.basketball-reference.com/players/m/miltosh01.html"

# Row 5: Tobias Harris
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = "Tobias Harris"
$ws.Range("D5").Value = "PF"
$ws.Range("E5").Value = "6-8"
$ws.Range("F5").Value = 226
$ws.Range("G5").Value = "July 15, 1992"
$ws.Range("H5").Value = "us"
$ws.Range("I5").Value = 11
$ws.Range("J5").Value = "Tennessee"
$ws.Range("K5").Value = "https://www.basketball-reference.com/players/h/harrito02.html"

# Row 6: P.J. Tucker
$ws.Range("B6").Value = 17
$ws.Range("C6").Value = "P.J. Tucker"
$ws.Range("D6").Value = "PF"
$ws.Range("E6").Value = "6-5"
$ws.Range("F6").Value = 245
$ws.Range("G6").Value = "May 5, 1985"
$ws.Range("H6").Value = "us"
$ws.Range("I6").Value = 11
$ws.Range("J6").Value = "Texas"
$ws.Range("K6").Value = "https://www.basketball-reference.com/players/t/tuckepj01.html"

# Row 9: James Harden
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "James Harden"
$ws.Range("D9").Value = "PG"
$ws.Range("E9").Value = "6-5"
$ws.Range("F9").Value = 220
$ws.Range("G9").Value = "August 26, 1989"
$ws.Range("H9").Value = "us"
$ws.Range("I9").Value = 13
$ws.Range("J9").Value = "Arizona State"
$ws.Range("K9").Value = "https://www.basketball-reference.com/players/h/hardeja01.html"

# Row 10: Tyrese Maxey
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = "Tyrese Maxey"
$ws.Range("D10").Value = "SG"
$ws.Range("E10").Value = "6-2"
$ws.Range("F10").Value = 200
$ws.Range("G10").Value = "November 4, 2000"
$ws.Range("H10").Value = "us"
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = "Kentucky"
$ws.Range("K10").Value = "https://www.basketball-reference.com/players/m/maxeyty01.html"

# Row 11: Montrezl Harrell
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = "Montrezl Harrell"
$ws.Range("D11").Value = "C"
$ws.Range("E11").Value = "6-7"
$ws.Range("F11").Value = 240
$ws.Range("G11").Value = "January 26, 1994"
$ws.Range("H11").Value = "us"
$ws.Range("I11").Value = 7
$ws.Range("J11").Value = "Louisville"
$ws.Range("K11").Value = "https://www.basketball-reference.com/players/h/harremo01.html"


Write-Output "Done applying roster reorder."
